$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A was empty; the edit shifts everything in B:E (values, styles and
# column widths) one column to the left into A:D. Deleting column A does
# exactly that.
$ws.Columns.Item(1).Delete()

# Deleting a column resets the custom width flag on the shifted columns in
# this engine, so restore the original widths (11, 5, 27, 12 characters)
# explicitly. COM's ColumnWidth is in character units, which the file
# format stores with a fixed +0.8333... padding offset, so subtract that
# offset to land back on the exact original widths on disk.
$offset = 0.8333333333333333
$ws.Columns.Item(1).ColumnWidth = 11 - $offset
$ws.Columns.Item(2).ColumnWidth = 5 - $offset
$ws.Columns.Item(3).ColumnWidth = 27 - $offset
$ws.Columns.Item(4).ColumnWidth = 12 - $offset
